$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force ambiguous numeric-looking strings to remain
# text without leaving a NumberFormat/style footprint on the target cell:
# format the scratch cell as Text, write the value there, copy it, and
# PasteSpecial -Values only onto the destination (values-only paste does not
# carry the "@" text format over, so the destination cell keeps its original
# (default) style while still inheriting the scratch cell's text storage).
$scratch = $ws.Range("ZZ1")

$ws.Range("D2").Value = '28.612.26'
$ws.Range("E2").Value = '  +0.65%  '

$ws.Range("D3").Value = '1.584.22'
$ws.Range("E3").Value = '  -0.44%  '

$scratch.NumberFormat = "@"
$scratch.Value = '213.55'
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
$ws.Range("E5").Value = '  +0.19%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.491'
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
$ws.Range("E6").Value = '  -0.20%  '

$ws.Range("E7").Value = '  -0.04%  '

$scratch.NumberFormat = "@"
$scratch.Value = '44.22'
$scratch.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
$ws.Range("E8").Value = '  -0.07%  '

$scratch.NumberFormat = "@"
$scratch.Value = '24.18'
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
$ws.Range("E9").Value = '  -0.97%  '

$ws.Range("E10").Value = '  -1.71%  '

$ws.Range("E11").Value = '  -1.31%  '

$ws.Range("E12").Value = '  +0.74%  '

$ws.Range("D13").Value = '1.809.98'
$ws.Range("E13").Value = '  -0.47%  '

$ws.Range("D14").Value = '1.584.25'
$ws.Range("E14").Value = '  -0.53%  '

$scratch.NumberFormat = "@"
$scratch.Value = '3.73'
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
$ws.Range("E15").Value = '  -0.51%  '

$ws.Range("D16").Value = '28.633.81'
$ws.Range("E16").Value = '  +0.61%  '

$ws.Range("E17").Value = '  -1.81%  '

$scratch.NumberFormat = "@"
$scratch.Value = '62.26'
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
$ws.Range("E18").Value = '  -1.36%  '

$scratch.NumberFormat = "@"
$scratch.Value = '231.99'
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
$ws.Range("E19").Value = '  +0.95%  '

$scratch.NumberFormat = "@"
$scratch.Value = '7.42'
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
$ws.Range("E20").Value = '  -0.72%  '

$ws.Range("E21").Value = '  -2.07%  '

$ws.Range("E22").Value = '  -0.07%  '

$ws.Range("E23").Value = '  -3.66%  '

$ws.Range("E24").Value = '  -1.53%  '

$ws.Range("E25").Value = '  +5.67%  '

$scratch.NumberFormat = "@"
$scratch.Value = '151.54'
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
$ws.Range("E26").Value = '  -0.14%  '

$ws.Range("E27").Value = '  -1.09%  '

$ws.Range("E28").Value = '  -1.65%  '

$ws.Range("E29").Value = '  -2.06%  '

$ws.Range("E30").Value = '  -0.08%  '

$ws.Range("E31").Value = '  +2.42%  '

$ws.Range("E32").Value = '  -1.67%  '

$ws.Range("E33").Value = '  -1.01%  '

$ws.Range("E34").Value = '  -1.71%  '

$ws.Range("D35").Value = '1.398.84'
$ws.Range("E35").Value = '  -0.18%  '

$ws.Range("E36").Value = '  +3.73%  '

$ws.Range("E37").Value = '  -3.82%  '

$ws.Range("E38").Value = '  +0.43%  '

$ws.Range("E39").Value = '  +3.50%  '

$ws.Range("E40").Value = '  -0.14%  '

$ws.Range("E41").Value = '  -3.44%  '

$ws.Range("E42").Value = '  -0.07%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.795'
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
$ws.Range("E43").Value = '  -2.06%  '

$ws.Range("E44").Value = '  +1.70%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.0465'
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
$ws.Range("E45").Value = '  -0.60%  '

$scratch.NumberFormat = "@"
$scratch.Value = '5.46'
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
$ws.Range("E46").Value = '  -2.72%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.960'
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
$ws.Range("E47").Value = '  -2.29%  '

$scratch.NumberFormat = "@"
$scratch.Value = '63.36'
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
$ws.Range("E48").Value = '  +0.16%  '

$ws.Range("D49").Value = '1.721.79'
$ws.Range("E49").Value = '  -0.27%  '

$scratch.NumberFormat = "@"
$scratch.Value = '86.85'
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
$ws.Range("E50").Value = '  -0.44%  '

$ws.Range("D51").Value = '0.0₆0103'
$ws.Range("E51").Value = '  -1.58%  '
